$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the F column values (budget allocations), powers of 10 instead of 1-8
$ws.Range("F4").Value = 10
$ws.Range("F5").Value = 100
$ws.Range("F6").Value = 1000
$ws.Range("F7").Value = 10000
$ws.Range("F8").Value = 100000
$ws.Range("F9").Value = 1000000
$ws.Range("F10").Value = 10000000
$ws.Range("F11").Value = 100000000

# Update selection to match the new active cell
$ws.Range("F20").Select()
